$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the most recent 20-row date block (1853:1872, date 44021) twice,
# to create rows for new dates 44022 and 44025, preserving styles/formatting.
$ws.Range("A1853:H1872").Copy($ws.Range("A1873:H1892"))
$ws.Range("A1853:H1872").Copy($ws.Range("A1893:H1912"))

# Update the date column and the numeric observation counts that differ
# from the template block for each of the two new dates.
$ws.Range("A1873").Value = 44022
$ws.Range("C1873").Value = 4
$ws.Range("D1873").Value = 1
$ws.Range("A1874").Value = 44022
$ws.Range("C1874").Value = 5
$ws.Range("A1875").Value = 44022
$ws.Range("A1876").Value = 44022
$ws.Range("C1876").Value = 1
$ws.Range("A1877").Value = 44022
$ws.Range("A1878").Value = 44022
$ws.Range("A1879").Value = 44022
$ws.Range("A1880").Value = 44022
$ws.Range("C1880").Value = 8
$ws.Range("D1880").Value = 1
$ws.Range("A1881").Value = 44022
$ws.Range("C1881").Value = 3
$ws.Range("A1882").Value = 44022
$ws.Range("A1883").Value = 44022
$ws.Range("A1884").Value = 44022
$ws.Range("A1885").Value = 44022
$ws.Range("A1886").Value = 44022
$ws.Range("A1887").Value = 44022
$ws.Range("A1888").Value = 44022
$ws.Range("A1889").Value = 44022
$ws.Range("A1890").Value = 44022
$ws.Range("A1891").Value = 44022
$ws.Range("A1892").Value = 44022
$ws.Range("A1893").Value = 44025
$ws.Range("C1893").Value = 10
$ws.Range("D1893").Value = ""
$ws.Range("A1894").Value = 44025
$ws.Range("C1894").Value = 15
$ws.Range("A1895").Value = 44025
$ws.Range("C1895").Value = ""
$ws.Range("A1896").Value = 44025
$ws.Range("C1896").Value = 1
$ws.Range("A1897").Value = 44025
$ws.Range("A1898").Value = 44025
$ws.Range("A1899").Value = 44025
$ws.Range("A1900").Value = 44025
$ws.Range("C1900").Value = 12
$ws.Range("A1901").Value = 44025
$ws.Range("C1901").Value = 3
$ws.Range("D1901").Value = 1
$ws.Range("A1902").Value = 44025
$ws.Range("A1903").Value = 44025
$ws.Range("A1904").Value = 44025
$ws.Range("C1904").Value = ""
$ws.Range("A1905").Value = 44025
$ws.Range("C1905").Value = 2
$ws.Range("A1906").Value = 44025
$ws.Range("A1907").Value = 44025
$ws.Range("A1908").Value = 44025
$ws.Range("C1908").Value = 1
$ws.Range("A1909").Value = 44025
$ws.Range("A1910").Value = 44025
$ws.Range("A1911").Value = 44025
$ws.Range("A1912").Value = 44025
